# Set all correlation values in B2:D9 to 0 (area check / specs files update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:D9").Value = 0
